$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns F, G, H (same header style as column E)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier flags (MAD-based) for rows 2-25
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false
$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false
$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $false
$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false
$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = $false
$ws.Range("F8").Value = $true
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false
$ws.Range("F9").Value = $false
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = $false
$ws.Range("F10").Value = $false
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = $false
$ws.Range("F11").Value = $false
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = $false
$ws.Range("F12").Value = $false
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = $false
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = $true
$ws.Range("H13").Value = $false
$ws.Range("F14").Value = $false
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = $false
$ws.Range("F15").Value = $false
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = $false
$ws.Range("F16").Value = $false
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = $false
$ws.Range("F17").Value = $false
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = $false
$ws.Range("F18").Value = $false
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = $false
$ws.Range("F19").Value = $false
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = $false
$ws.Range("F20").Value = $false
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = $false
$ws.Range("F21").Value = $false
$ws.Range("G21").Value = $false
$ws.Range("H21").Value = $false
$ws.Range("F22").Value = $true
$ws.Range("G22").Value = $true
$ws.Range("H22").Value = $false
$ws.Range("F23").Value = $false
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = $false
$ws.Range("F24").Value = $false
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = $false
$ws.Range("F25").Value = $false
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = $false
